# Auto-generated edit script: updates computed market-price columns (H-N)
# across multiple sheets to match the refreshed source data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 801.6667
$ws.Range("I28").Value = 510.36365
$ws.Range("K28").Value = 510.36365
$ws.Range("M28").Value = -25.36365000000001
$ws.Range("H58").Value = 1889.8462
$ws.Range("I58").Value = 946.3333
$ws.Range("J58").Value = 4012.75
$ws.Range("K58").Value = 2838.9999
$ws.Range("L58").Value = 12038.25
$ws.Range("M58").Value = -2688.9999
$ws.Range("N58").Value = -12338.25
$ws.Range("H132").Value = 2605.3333
$ws.Range("I132").Value = 1550.8334
$ws.Range("K132").Value = 4652.5002
$ws.Range("M132").Value = -2122.5002
$ws.Range("H135").Value = 620.14813
$ws.Range("I135").Value = 455.53845
$ws.Range("K135").Value = 4099.84605
$ws.Range("M135").Value = -1564.84605

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2058.8262
$ws.Range("I32").Value = 2077.975
$ws.Range("J32").Value = 1931.1666
$ws.Range("K32").Value = 2077.975
$ws.Range("L32").Value = 1931.1666
$ws.Range("M32").Value = -1790.975
$ws.Range("N32").Value = -2505.1666
$ws.Range("H44").Value = 22198.1
$ws.Range("J44").Value = 22198.1
$ws.Range("L44").Value = 22198.1
$ws.Range("N44").Value = -23174.1
$ws.Range("H46").Value = 12000
$ws.Range("I46").Value = 12000
$ws.Range("J46").Value = 12000
$ws.Range("K46").Value = 12000
$ws.Range("L46").Value = 12000
$ws.Range("M46").Value = -11681
$ws.Range("N46").Value = -12638
$ws.Range("H55").Value = 22252.7
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 22252.7
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 22252.7
$ws.Range("M55").Value = $null
$ws.Range("N55").Value = -22882.7
$ws.Range("H74").Value = 4631639.5
$ws.Range("I74").Value = 2647674
$ws.Range("J74").Value = 9260892
$ws.Range("K74").Value = 2647674
$ws.Range("L74").Value = 9260892
$ws.Range("M74").Value = -2646800
$ws.Range("N74").Value = -9262640
$ws.Range("H77").Value = 4631639.5
$ws.Range("I77").Value = 2647674
$ws.Range("J77").Value = 9260892
$ws.Range("K77").Value = 13238370
$ws.Range("L77").Value = 46304460
$ws.Range("M77").Value = -13234002
$ws.Range("N77").Value = -46313196
$ws.Range("H80").Value = 20805
$ws.Range("H83").Value = 20805
$ws.Range("H97").Value = 1173.8572
$ws.Range("I97").Value = 804.5
$ws.Range("K97").Value = 804.5
$ws.Range("M97").Value = -308.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12224.75
$ws.Range("I94").Value = 2949.5
$ws.Range("K94").Value = 2949.5
$ws.Range("M94").Value = -2498.5
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = $null
$ws.Range("H105").Value = 2863
$ws.Range("I105").Value = 1989
$ws.Range("J105").Value = 3300
$ws.Range("K105").Value = 1989
$ws.Range("L105").Value = 3300
$ws.Range("M105").Value = -242
$ws.Range("N105").Value = -6794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = $null
$ws.Range("H99").Value = 2103.9375
$ws.Range("I99").Value = 2000.7
$ws.Range("J99").Value = 2276
$ws.Range("K99").Value = 2000.7
$ws.Range("L99").Value = 2276
$ws.Range("M99").Value = -502.7
$ws.Range("N99").Value = -5272
$ws.Range("H126").Value = 2103.9375
$ws.Range("I126").Value = 2000.7
$ws.Range("J126").Value = 2276
$ws.Range("K126").Value = 6002.1
$ws.Range("L126").Value = 6828
$ws.Range("M126").Value = -3532.1
$ws.Range("N126").Value = -11768
$ws.Range("H134").Value = 2780703
$ws.Range("I134").Value = 2229
$ws.Range("K134").Value = 6687
$ws.Range("M134").Value = -4152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 264
$ws.Range("I92").Value = 269.54544
$ws.Range("J92").Value = 203
$ws.Range("K92").Value = 808.63632
$ws.Range("L92").Value = 609
$ws.Range("M92").Value = 439.36368
$ws.Range("N92").Value = -3105
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3439.6956
$ws.Range("I122").Value = 3902.9285
$ws.Range("K122").Value = 11708.7855
$ws.Range("M122").Value = -9258.7855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4032.6667
$ws.Range("I7").Value = 3824
$ws.Range("J7").Value = 4450
$ws.Range("K7").Value = 3824
$ws.Range("L7").Value = 4450
$ws.Range("M7").Value = -3712
$ws.Range("N7").Value = -4674
$ws.Range("H40").Value = 2199.1538
$ws.Range("I40").Value = 2244.4546
$ws.Range("K40").Value = 2244.4546
$ws.Range("M40").Value = -2108.4546
$ws.Range("H68").Value = 1685.5
$ws.Range("J68").Value = 1772.75
$ws.Range("L68").Value = 1772.75
$ws.Range("N68").Value = -3270.75
$ws.Range("H71").Value = 1685.5
$ws.Range("J71").Value = 1772.75
$ws.Range("L71").Value = 8863.75
$ws.Range("N71").Value = -16351.75
$ws.Range("H82").Value = 1624.4
$ws.Range("I82").Value = 680.25
$ws.Range("K82").Value = 680.25
$ws.Range("M82").Value = -319.25
$ws.Range("H85").Value = 1624.4
$ws.Range("I85").Value = 680.25
$ws.Range("K85").Value = 680.25
$ws.Range("M85").Value = 567.75
$ws.Range("H93").Value = 2666.5557
$ws.Range("I93").Value = 2750
$ws.Range("K93").Value = 2750
$ws.Range("M93").Value = -1502
$ws.Range("H126").Value = 4032.6667
$ws.Range("I126").Value = 3824
$ws.Range("J126").Value = 4450
$ws.Range("K126").Value = 11472
$ws.Range("L126").Value = 13350
$ws.Range("M126").Value = -9002
$ws.Range("N126").Value = -18290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 22183.666
$ws.Range("I45").Value = 20477.166
$ws.Range("K45").Value = 20477.166
$ws.Range("M45").Value = -19986.166
$ws.Range("H61").Value = 7624
$ws.Range("I61").Value = 1539
$ws.Range("J61").Value = 9652.333000000001
$ws.Range("K61").Value = 1539
$ws.Range("L61").Value = 9652.333000000001
$ws.Range("M61").Value = -1247
$ws.Range("N61").Value = -10236.333
